$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Update "Förändrad" date (column C) for rows 2 through 15 from 45185 (2023-09-16)
# to 45204 (2023-10-05), as per the diff.
for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45204
}
